# "braun and miele added"
#
# The sheet originally had AEG / Illy / (Illy note) in rows 9-11 and
# SMEG in row 17 (with a couple of leftover/blank helper rows in
# between). The edit inserts a new "Miele" row and a "Braun" block
# (manufacturer row + a bunch of plain-text manual links) right after
# Jura (row 8), which pushes SMEG up to row 22 and AEG / Illy / the Illy
# note down to rows 23-25.
#
# We rebuild the final layout directly rather than trying to replay the
# exact cut/paste the author did in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Move the existing AEG / Illy / Illy-note / SMEG rows down/up to
#    their new homes (rows 23, 24, 25, 22) before anything else so we
#    don't clobber data we still need to read.
# ---------------------------------------------------------------

# SMEG: old row 17 -> new row 22
$ws.Range("A22").Value = $ws.Range("A17").Value()
$ws.Range("B22").Value = $ws.Range("B17").Value()
$ws.Range("B22").Style = "Link"
$ws.Range("C22").Value = $ws.Range("C17").Value()

# AEG: old row 9 -> new row 23
$ws.Range("A23").Value = $ws.Range("A9").Value()
$ws.Range("B23").Value = $ws.Range("B9").Value()
$ws.Range("B23").Style = "Link"
$ws.Range("C23").Value = $ws.Range("C9").Value()
$ws.Range("D23").Value = $ws.Range("D9").Value()

# Illy: old row 10 -> new row 24
$ws.Range("A24").Value = $ws.Range("A10").Value()
$ws.Range("B24").Value = $ws.Range("B10").Value()
$ws.Range("B24").Style = "Link"
$ws.Range("C24").Value = $ws.Range("C10").Value()
$ws.Range("D24").Value = $ws.Range("D10").Value()
$ws.Range("K24").Value = $ws.Range("K10").Value()

# Illy note: old row 11 -> new row 25
$ws.Range("B25").Value = $ws.Range("B11").Value()

# ---------------------------------------------------------------
# 2) Clear out the old locations that are not part of the new block
#    (row 17 and the stray notes in rows 13/16 get overwritten below
#    anyway; row 10/11 get overwritten by the new Braun rows below).
# ---------------------------------------------------------------
$ws.Range("A9:K11").ClearContents()
$ws.Range("A17:K17").ClearContents()

# ---------------------------------------------------------------
# 3) Write the new Miele row (row 9).
# ---------------------------------------------------------------
$ws.Range("A9").Value = "Miele"
$ws.Range("B9").Value = "https://www.miele.de/haushalt/produktauswahl-kaffeevollautomaten-2513.htm"
$ws.Range("C9").Value = "                                    "

# ---------------------------------------------------------------
# 4) Write the new Braun block (rows 10-17).
# ---------------------------------------------------------------
$ws.Range("A10").Value = "Braun"
$ws.Range("B10").Value = "https://www.braunhousehold.com/de-de/produkte/fruehstueck/kaffeemaschinen/c/coffee_machines?q=%3Arelevance%3Acategory_cluster_coffee_makers%3AModerne%2BKaffeemaschine%3Acategory_cluster_coffee_makers%3AStandard-Kaffeemaschine"

# Row 11: manuals header ("bedienungsanleitungen ->"), stored with a
# leading apostrophe so Excel keeps it as quote-prefixed text (this is
# what produces the extra quotePrefix cell style in the saved file).
$ws.Range("A11").Value = "'bedienungsanleitungen ->"
$ws.Range("B11").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/aromaster-kaffeemaschinen/c/aromaster_coffee_machines"

$ws.Range("B12").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/puraroma-7/c/puraroma_7"

$ws.Range("B13").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/id-fruehstuecksserie-kaffeemaschinen/c/id_breakfast_collection_coffee_machines"

$ws.Range("B14").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/purease-kaffeemaschinen/c/purease_coffee_machines"
$ws.Range("C14").Value = "                                    "

$ws.Range("B15").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/cafehouse-coffee-machines/c/cafehouse_coffee_machines"

$ws.Range("B16").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/multiserve-kaffeemaschine/c/multiserve_coffee_machines"
$ws.Range("B16").Style = "Link"

$ws.Range("B17").Value = "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/breakfast-1-collection-coffee-machines/c/breakfast_1_collection_coffee_machines"

# Only B14 actually becomes a real hyperlink (the rest of the Braun
# manual links stay plain text, matching the source workbook).
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.braunhousehold.com/de-de/manuals/produkte/fruehstueck/kaffeemaschinen/purease-kaffeemaschinen/c/purease_coffee_machines", [Type]::Missing, [Type]::Missing, $ws.Range("B14").Value())
$ws.Range("B14").Style = "Link"

# ---------------------------------------------------------------
# 5) Re-point the window/selection the way the saved file shows
#    (scrolled so row 10 is at the top, with B11:B17 selected).
# ---------------------------------------------------------------
$ws.Range("B11:B17").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
